$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 36

$ws.Cells.Item($row, 1).Value = "'03-11-2021"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 10000
$ws.Cells.Item($row, 4).Value = 0
